# AoC 2024 - Excel update
# Adds a new "Day 20" worksheet (after "Day 17") containing a small
# racetrack/maze grid, and makes it the active/selected sheet.

$wb = $excel.ActiveWorkbook

# Add the new sheet at the end of the tab strip (after the current last sheet).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Day 20"

# Fill in the little maze grid. Write the start tile ("S") first, then the
# wall tiles ("#"), then the open tiles ("."), matching the order the values
# were originally typed in (and therefore the order they land in the shared
# string table).
$ws.Range("D5").Value = "S"

$ws.Range("D4").Value = "#"
$ws.Range("C5").Value = "#"
$ws.Range("E5").Value = "#"
$ws.Range("D6").Value = "#"

$ws.Range("D3").Value = "."
$ws.Range("C4").Value = "."
$ws.Range("E4").Value = "."
$ws.Range("B5").Value = "."
$ws.Range("F5").Value = "."
$ws.Range("C6").Value = "."
$ws.Range("E6").Value = "."
$ws.Range("D7").Value = "."

# Make "Day 20" the active/visible tab with its own saved view state.
$ws.Activate()
$excel.ActiveWindow.Zoom = 175
$ws.Range("J13").Select() | Out-Null
